$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price entry is added at the top of the data (row 2). Every
# existing data row (2-14) shifts down by one (3-15). Insert the new row
# just below the current row 2 (not above it) so the freshly-inserted row
# inherits row 2's plain data formatting instead of row 1's bold header
# formatting; then re-home the old row-2 values into the new row 3 and
# write the new entry's values into row 2.
$ws.Rows.Item(3).Insert()

# Row 3: the data that used to live in row 2 (unchanged values, just moved down)
$ws.Range("A3").Value = 10
$ws.Range("B3").Value = "Vega Modelo de Temuco"
$ws.Range("C3").Value = "La Araucanía"
$ws.Range("D3").Value = 44418
$ws.Range("E3").Value = 9
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100104
$ws.Range("H3").Value = "Frutos de pepita"
$ws.Range("I3").Value = 100104005
$ws.Range("J3").Value = "Pera asiática"
$ws.Range("K3").Value = "Hosui"
$ws.Range("L3").Value = "Especial"
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 8000
$ws.Range("O3").Value = 8000
$ws.Range("P3").Value = 8000
$ws.Range("Q3").Value = "$/caja 15 kilos granel"
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 533
$ws.Range("T3").Value = 15

# Row 2: the new weekly entry
$ws.Range("A2").Value = 10
$ws.Range("B2").Value = "Vega Modelo de Temuco"
$ws.Range("C2").Value = "La Araucanía"
$ws.Range("D2").Value = 45083
$ws.Range("E2").Value = 9
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100104
$ws.Range("H2").Value = "Frutos de pepita"
$ws.Range("I2").Value = 100104005
$ws.Range("J2").Value = "Pera asiática"
$ws.Range("K2").Value = "Hosui"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 55
$ws.Range("N2").Value = 16000
$ws.Range("O2").Value = 16000
$ws.Range("P2").Value = 16000
$ws.Range("Q2").Value = "$/caja 18 kilos granel"
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 889
$ws.Range("T2").Value = 18
